$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9152
$ws1.Range("F13").Value = 9371
$ws1.Range("F16").Value = 216
$ws1.Range("F22").Value = 132
$ws1.Range("F26").Value = 52
$ws1.Range("F27").Value = 93
$ws1.Range("F30").Value = 105
$ws1.Range("F33").Value = 907
$ws1.Range("F38").Value = 399
$ws1.Range("F39").Value = 197
$ws1.Range("F42").Value = 162
$ws1.Range("F44").Value = 45
$ws1.Range("F45").Value = 94
$ws1.Range("F46").Value = 57
$ws1.Range("F47").Value = 20
$ws1.Range("F48").Value = 4024
$ws1.Range("F49").Value = 31

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 22
$ws2.Range("F18").Value = 40

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 374

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 9152
$ws4.Range("F5").Value = 374
$ws4.Range("F16").Value = 9371
$ws4.Range("F19").Value = 216
$ws4.Range("F26").Value = 52
$ws4.Range("F35").Value = 399
$ws4.Range("F37").Value = 197
$ws4.Range("F41").Value = 45
$ws4.Range("F42").Value = 94
$ws4.Range("F43").Value = 57
$ws4.Range("F45").Value = 4024
